$wb = $excel.ActiveWorkbook

# Helper: duplicate formatting from a template row ($templateRange, e.g. "A4:N4")
# onto the target row, then write the given values across columns A..N.
function Add-ScrimRow {
    param($ws, $row, $templateRange, $values)
    $ws.Range($templateRange).Copy()
    $ws.Range("A$row`:N$row").PasteSpecial(-4122)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# ---- Sheet: Triple Dribble ----
$ws = $wb.Worksheets.Item("Triple Dribble")

Add-ScrimRow $ws 57 "A4:N4" @("HANK","BEA","LUMI","CORDELIUS","CHARLIE","TICK","Equipo 2","PLP|BrriN","MTM|snoiy","PLP|Mine","NHG|Xemp","NHG|Bayarea","NHG|GN","20250724T004228.000Z")
Add-ScrimRow $ws 58 "A4:N4" @("HANK","BEA","LUMI","CORDELIUS","CHARLIE","TICK","Equipo 2","PLP|BrriN","MTM|snoiy","PLP|Mine","NHG|Xemp","NHG|Bayarea","NHG|GN","20250724T004029.000Z")
Add-ScrimRow $ws 59 "A5:N5" @("HANK","BEA","LUMI","CORDELIUS","CHARLIE","TICK","Equipo 1","PLP|BrriN","MTM|snoiy","PLP|Mine","NHG|Xemp","NHG|Bayarea","NHG|GN","20250724T003739.000Z")
Add-ScrimRow $ws 60 "A5:N5" @("JAE-YONG","MOE","ASH","HANK","MEEPLE","SHADE","Equipo 1","MTM|snoiy","PLP|Mine","PLP|BrriN","NHG|Xemp","NHG|Bayarea","NHG|GN","20250724T003149.000Z")
Add-ScrimRow $ws 61 "A5:N5" @("JAE-YONG","MOE","ASH","HANK","MEEPLE","SHADE","Equipo 1","MTM|snoiy","PLP|Mine","PLP|BrriN","NHG|Xemp","NHG|Bayarea","NHG|GN","20250724T002859.000Z")
$excel.CutCopyMode = $false

# ---- Sheet: Pinball Dreams ----
$ws = $wb.Worksheets.Item("Pinball Dreams")

Add-ScrimRow $ws 9 "A4:N4" @("STU","BULL","BEA","BIBI","MANDY","CORDELIUS","Equipo 1","Finki is back.","Solar Ray ☀️","Xyz","BC*|Jubileubr","LOUD|Edinho","CASA|Doritos","20250724T002248.000Z")
Add-ScrimRow $ws 10 "A4:N4" @("STU","BULL","BEA","BIBI","MANDY","CORDELIUS","Equipo 1","Finki is back.","Solar Ray ☀️","Xyz","BC*|Jubileubr","LOUD|Edinho","CASA|Doritos","20250724T001957.000Z")
$excel.CutCopyMode = $false

# ---- Sheet: Crystal Arcade ----
$ws = $wb.Worksheets.Item("Crystal Arcade")

Add-ScrimRow $ws 36 "A4:N4" @("DRACO","BARLEY","MEG","MEEPLE","EMZ","CORDELIUS","Equipo 1","TE|Rafikii","TE|Ezlivi","TE|Belal","TRB|Zeus 解開","TRB|R B M","TRB|Lxffy","20250724T004041.000Z")
Add-ScrimRow $ws 37 "A6:N6" @("DRACO","MEG","BERRY","MEEPLE","CORDELIUS","ASH","Equipo 2","TE|Rafikii","TE|Ezlivi","TE|Belal","TRB|Zeus 解開","TRB|Lxffy","TRB|R B M","20250724T003715.000Z")
Add-ScrimRow $ws 38 "A6:N6" @("DRACO","MEG","BERRY","MEEPLE","CORDELIUS","ASH","Equipo 2","TE|Rafikii","TE|Ezlivi","TE|Belal","TRB|Zeus 解開","TRB|Lxffy","TRB|R B M","20250724T003540.000Z")
Add-ScrimRow $ws 39 "A4:N4" @("CROW","LOU","KAZE","JESSIE","FINX","DOUG","Equipo 1","TE|Rafikii","TE|Ezlivi","TE|Belal","TRB|Zeus 解開","TRB|Lxffy","TRB|R B M","20250724T002751.000Z")
Add-ScrimRow $ws 40 "A4:N4" @("CROW","LOU","KAZE","JESSIE","FINX","DOUG","Equipo 1","TE|Rafikii","TE|Ezlivi","TE|Belal","TRB|Zeus 解開","TRB|Lxffy","TRB|R B M","20250724T002523.000Z")
Add-ScrimRow $ws 41 "A6:N6" @("CROW","LOU","KAZE","JESSIE","FINX","DOUG","Equipo 2","TE|Rafikii","TE|Ezlivi","TE|Belal","TRB|Zeus 解開","TRB|Lxffy","TRB|R B M","20250724T002308.000Z")
Add-ScrimRow $ws 42 "A6:N6" @("AMBER","BUZZ","ASH","SAM","STU","KENJI","Equipo 2","Finki is back.","Solar Ray ☀️","Xyz","BC*|Jubileubr","LOUD|Edinho","LOUD|KaioDog","20250724T003147.000Z")
Add-ScrimRow $ws 43 "A6:N6" @("AMBER","BUZZ","ASH","SAM","STU","KENJI","Equipo 2","Finki is back.","Solar Ray ☀️","Xyz","BC*|Jubileubr","LOUD|Edinho","LOUD|KaioDog","20250724T002947.000Z")
$excel.CutCopyMode = $false

# ---- Sheet: Layer Cake ----
$ws = $wb.Worksheets.Item("Layer Cake")

Add-ScrimRow $ws 48 "A4:N4" @("DOUG","PENNY","JANET","GUS","MR. P","CORDELIUS","Equipo 2","TE|Rafikii","TE|Ezlivi","TE|Belal","TRB|Zeus 解開","TRB|R B M","TRB|Lxffy","20250724T004857.000Z")
Add-ScrimRow $ws 49 "A8:N8" @("DOUG","PENNY","JANET","GUS","MR. P","CORDELIUS","Equipo 1","TE|Rafikii","TE|Ezlivi","TE|Belal","TRB|Zeus 解開","TRB|R B M","TRB|Lxffy","20250724T004637.000Z")
$excel.CutCopyMode = $false

# ---- Sheet: Dry Season ----
$ws = $wb.Worksheets.Item("Dry Season")

Add-ScrimRow $ws 36 "A4:N4" @("BONNIE","LUMI","HANK","MANDY","KAZE","BELLE","Equipo 1","Finki is back.","Solar Ray ☀️","Xyz","LOUD|Edinho","CASA|Doritos","LOUD|KaioDog","20250724T004000.000Z")
Add-ScrimRow $ws 37 "A4:N4" @("BONNIE","LUMI","HANK","MANDY","KAZE","BELLE","Equipo 1","Finki is back.","Solar Ray ☀️","Xyz","LOUD|Edinho","CASA|Doritos","LOUD|KaioDog","20250724T003739.000Z")
$excel.CutCopyMode = $false
